$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("The histogram")
$co = $ws.Shapes.AddChart2(-1, 118)
